# Apply numeric updates to the "Valores" sheet of the workbook that feeds
# the Power BI report (backup-generation helper described in the commit
# message). Only cell values change - no structural edits are required.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Valores")

$ws.Range("C1").Value = 5
$ws.Range("D1").Value = 2
$ws.Range("E1").Value = 1
$ws.Range("F1").Value = 2
$ws.Range("I1").Value = 0
$ws.Range("J1").Value = 3
$ws.Range("K1").Value = 1

$ws.Range("G2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0

$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 1

$ws.Range("C5").Value = 1

$ws.Range("C6").Value = 2
$ws.Range("D6").Value = 1
$ws.Range("F6").Value = 1
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 2
$ws.Range("K6").Value = 1

$ws.Range("J10").Value = 0

$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0

$ws.Range("C16").Value = 9
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 3
$ws.Range("G16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0
